$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Range("A1").Value = "test"
